# "Turklerin temmuz avansini hazirana gecirdim." - move the July advance
# payments forward into June's AVANS (H) column for the affected workers,
# which reduces their NET HAKEDIS (I) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# AVANS (H column) updates, row by row:
$ws.Range("H3").Value  = 10000   # SADIK ACAR:      0 -> 10000
$ws.Range("H4").Value  = 15000   # IHSAN GOL:    5000 -> 15000
$ws.Range("H5").Value  = 5000    # YUSUF TOMAK:     0 -> 5000
$ws.Range("H6").Value  = 5000    # KADIR MISIRLI:   0 -> 5000
$ws.Range("H8").Value  = 5000    # YUSUF AKKOYUN:   0 -> 5000
$ws.Range("H9").Value  = 5000    # GURKAN AKTAS:    0 -> 5000
$ws.Range("H10").Value = 13000   # BEKIR KOCAK:  3000 -> 13000
